# horario.xlsx - actualizamos horario y agregamos el primer trabajo de computacion numerica
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Color constants (BGR integers as used by Excel COM Interior.Color / Font.Color)
$red        = 255        # FF0000
$blue       = 15773696   # 00B0F0
$green      = 5287936    # 00B050
$lightgreen = 5296274    # 92D050
$darkblue   = 12611584   # 0070C0
$purple     = 10498160   # 7030A0

$xlCenter = -4108

function Clear-Cell($addr) {
    $r = $ws.Range($addr)
    $r.ClearContents()
    $r.Style = "Normal"
    return $r
}

function Set-CellText($addr, $text, $fill, $wrap) {
    $r = Clear-Cell $addr
    if ($text -ne $null -and $text -ne "") {
        $r.Value = $text
    }
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
    if ($wrap) {
        $r.WrapText = $true
    } else {
        $r.WrapText = $false
    }
    if ($fill -ne $null) {
        $r.Interior.Color = $fill
    }
    return $r
}

# --- Row 2 (8:30 - 9:30) ---
# F2: "Sistema operativo" -> now blank (keep centered/wrapped, no fill)
Set-CellText "F2" $null $null $true | Out-Null

# --- Row 3 (9:35 - 10:35) ---
# F3: "Sistema operativo" -> now blank, with red font color left as a marker
$r = Set-CellText "F3" $null $null $true
$r.Font.Color = $red

# --- Row 4 (10:50 - 11:50) ---
# B4: "Etica cristiana" -> blank (no fill, center/vcenter, no wrap)
Set-CellText "B4" $null $null $false | Out-Null
# E4: "Logica" -> blank (no fill, center/vcenter, wrap)
Set-CellText "E4" $null $null $true | Out-Null
# F4: "Sistema operativo" -> blank (no fill, center/vcenter, wrap)
Set-CellText "F4" $null $null $true | Out-Null

# --- Row 5 (11:55 - 12:55) ---
# B5: "Etica cristiana" -> blank
Set-CellText "B5" $null $null $false | Out-Null
# C5: "Base de datos Lab" -> blank
Set-CellText "C5" $null $null $true | Out-Null
# E5: "Logica" -> blank
Set-CellText "E5" $null $null $true | Out-Null
# F5: "Logica Lab" -> "Sistema operativo" (blue fill)
Set-CellText "F5" ("Sistema" + [char]10 + "operativo") $blue $true | Out-Null
# H5: keep underlined marker cell (style renumbers only, same visuals)
$r = Clear-Cell "H5"
$r.Font.Underline = 2

# --- Row 6 (13:10 - 14:10) ---
# C6: "Base de datos Lab" -> blank
Set-CellText "C6" $null $null $true | Out-Null
# F6: blank -> blank (no visual change, just touch formatting)
Set-CellText "F6" $null $null $false | Out-Null

# --- Row 7 (14:30 - 15:30) ---
# E7: "Base de datos" -> "Base de datos Lab" (green fill, now wrapped)
Set-CellText "E7" ("Base de datos" + [char]10 + "Lab") $green $true | Out-Null
# F7: blank -> "Sistema operativo" (blue fill)
Set-CellText "F7" ("Sistema" + [char]10 + "operativo") $blue $true | Out-Null

# --- Row 8 (15:35 - 16:35) ---
# E8: "Base de datos" -> "Base de datos Lab" (green fill, now wrapped)
Set-CellText "E8" ("Base de datos" + [char]10 + "Lab") $green $true | Out-Null
# F8: blank -> "Sistema operativo" (blue fill)
Set-CellText "F8" ("Sistema" + [char]10 + "operativo") $blue $true | Out-Null

# --- Row 9 (16:50 - 17:50) ---
# D9: "Programacion avanzada" -> blank
Set-CellText "D9" $null $null $true | Out-Null
# F9: "Programacion avanzada" -> blank
Set-CellText "F9" $null $null $true | Out-Null

# --- Row 10 (17:55 - 18:55) ---
# D10: "Programacion avanzada" -> "Base de datos" (green fill)
Set-CellText "D10" "Base de datos" $green $true | Out-Null
# F10: "Programacion avanzada" -> blank
Set-CellText "F10" $null $null $true | Out-Null

# --- Row 11 (19:10 - 20:10) ---
# D11: blank -> "Base de datos" (green fill, center/vcenter, no wrap)
Set-CellText "D11" "Base de datos" $green $false | Out-Null

# Update the active selection to match the saved workbook state (L4)
$ws.Range("L4").Select() | Out-Null
